# worked on inventory bugs solve
#
# - Remove the "Location" column (column K) entirely (header + data).
# - Rename row-2 sub-category/product values:
#     "AAC BLOCK" -> "BLOCKS"
#     "AAC FLY ASH BLOCKS(600*200*150) 2" -> "SOFT BLOCKS"
# - HSN/SAC Code for row 2 becomes the numeric value 638271023
#   (was the text "68159910").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column K (Location header + TAMARA value) - shifts dimension to A1:J3
$ws.Range("K1:K2").EntireColumn.Delete()

# Update row 2 (CIVIL MATERIAL / ... / CBM / 68159910 / material / 3650)
$ws.Range("B2").Value = "BLOCKS"
$ws.Range("C2").Value = "SOFT BLOCKS"
$ws.Range("E2").Value = 638271023

# Keep the active selection consistent with the edited range
$ws.Range("A3:G3").Select()
